$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value that must be stored as TEXT (shared string) even
# though it "looks like" a number or a date (e.g. "51336" or "2026-01-30").
# Plain .Value assignment would let the engine infer a numeric/date type and
# stamp a new number-format style on the cell; forcing the cell to text
# format first (then resetting the style afterwards) keeps the underlying
# cell style index unchanged while still yielding a text cell.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ===========================================================================
# Sheet "Reports" - append two new report rows (34 and 35)
# ===========================================================================
$wsReports = $wb.Worksheets.Item("Reports")
$loReports = $wsReports.ListObjects.Item(1)

$loReports.ListRows.Add() | Out-Null
$loReports.ListRows.Add() | Out-Null

# Row 34
Set-TextValue $wsReports.Cells.Item(34, 1) "51336"
$wsReports.Cells.Item(34, 2).Value = "HC 1658"
$wsReports.Cells.Item(34, 3).Value = "2024-26"
$wsReports.Cells.Item(34, 4).Value = "Defence Committee"
$wsReports.Cells.Item(34, 5).Value = "Commons"
$wsReports.Cells.Item(34, 6).Value = "The UK contribution to European Security: Government Response"
$wsReports.Cells.Item(34, 7).Value = "6th Special Report"
Set-TextValue $wsReports.Cells.Item(34, 8) "2026-01-30"
$wsReports.Cells.Item(34, 9).Value = "11:00:00"
$wsReports.Cells.Item(34, 11).Value = "0:55:37"

# Row 35
Set-TextValue $wsReports.Cells.Item(35, 1) "51345"
$wsReports.Cells.Item(35, 2).Value = "HC 291-xlvi"
$wsReports.Cells.Item(35, 3).Value = "2024-26"
$wsReports.Cells.Item(35, 4).Value = "Statutory Instruments (Joint Committee)"
$wsReports.Cells.Item(35, 5).Value = "Joint"
$wsReports.Cells.Item(35, 6).Value = "Forty-sixth Report - 3 Statutory Instruments Reported"
Set-TextValue $wsReports.Cells.Item(35, 8) "2026-01-30"
$wsReports.Cells.Item(35, 9).Value = "11:00:00"
$wsReports.Cells.Item(35, 11).Value = "0:55:37"

# ===========================================================================
# Sheet "Scans" - append one new scan row (21)
# ===========================================================================
$wsScans = $wb.Worksheets.Item("Scans")
$loScans = $wsScans.ListObjects.Item(1)

$loScans.ListRows.Add() | Out-Null

Set-TextValue $wsScans.Cells.Item(21, 1) "2026-01-30"
$wsScans.Cells.Item(21, 2).Value = "11:55:37"
$wsScans.Cells.Item(21, 3).Value = "51336, 51345"
